$d = $word.ActiveDocument

# Remember how many paragraphs exist before the edit; the last one is the
# "Sim, SCRUM poderia..." paragraph that the new content is appended after.
$baseCount = $d.Paragraphs.Count

# Create all the new (still-empty) paragraphs first, *before* any direct
# character formatting is applied anywhere, so nothing can inherit bold /
# font-size from a neighbouring run.
$end = $d.Paragraphs.Last.Range
$end.Collapse(0)
$end.InsertParagraphAfter()        # baseCount+1 : blank separator paragraph

$end = $d.Paragraphs.Last.Range
$end.Collapse(0)
$end.InsertParagraphAfter()        # baseCount+2 : "Questão 4)" heading

$end = $d.Paragraphs.Last.Range
$end.Collapse(0)
$end.InsertParagraphAfter()        # baseCount+3

$end = $d.Paragraphs.Last.Range
$end.Collapse(0)
$end.InsertParagraphAfter()        # baseCount+4

$end = $d.Paragraphs.Last.Range
$end.Collapse(0)
$end.InsertParagraphAfter()        # baseCount+5

$end = $d.Paragraphs.Last.Range
$end.Collapse(0)
$end.InsertParagraphAfter()        # baseCount+6

$end = $d.Paragraphs.Last.Range
$end.Collapse(0)
$end.InsertParagraphAfter()        # baseCount+7

$end = $d.Paragraphs.Last.Range
$end.Collapse(0)
$end.InsertParagraphAfter()        # baseCount+8

# Fill in the text of each new paragraph (still unformatted / plain).
$d.Paragraphs.Item($baseCount + 2).Range.Text = "Questão 4)"
$d.Paragraphs.Item($baseCount + 3).Range.Text = "Começaria organizando a equipe em três equipes:"
$d.Paragraphs.Item($baseCount + 4).Range.Text = "Equipe do Estado do trânsito na rota;"
$d.Paragraphs.Item($baseCount + 5).Range.Text = "Equipe da Geolocalização;"
$d.Paragraphs.Item($baseCount + 6).Range.Text = "Equipe do Apontamento de lixo coletado;"
$d.Paragraphs.Item($baseCount + 7).Range.Text = "Poderia ocorrer mudanças caso o cliente pedir, ou a equipe implementar, uma outra funcionalidade."
$d.Paragraphs.Item($baseCount + 8).Range.Text = "Dentro deles teriam funções como Administrador de Dados, organizar dados sobre o lixo coletado e ruas onde tem quantidade maiores de lixo, Analista de Processos e Negócio, Arquiteto de Solução, Engenheiro de Software e entre outros."

# Finally, make the "Questão 4)" paragraph a bold, 12pt heading, matching
# the style used for "Questão 2)" / "Questão 3)" above it. This is done
# last so the bold/size never bleeds forward into the later paragraphs.
$headingRange = $d.Paragraphs.Item($baseCount + 2).Range
$headingRange.Font.Bold = $true
$headingRange.Font.BoldBi = $true
$headingRange.Font.Size = 12
$headingRange.Font.SizeBi = 12

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
